# Update Name of Algo
# Applies updated RandomForest imputation values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.3834
$ws.Range("B7").Value = 5.6021
$ws.Range("A9").Value = -21.73280000000001
$ws.Range("B12").Value = 5.833199999999994
$ws.Range("B14").Value = 6.133700000000005
$ws.Range("E15").Value = 16.46899999999999
$ws.Range("A18").Value = -22.28870000000002
$ws.Range("A20").Value = -19.53909999999999
$ws.Range("B26").Value = 3.991300000000003
$ws.Range("A27").Value = -22.06809999999999
$ws.Range("B27").Value = 5.397200000000003
$ws.Range("B29").Value = 4.9015
$ws.Range("E33").Value = 17.11500000000002
$ws.Range("A35").Value = -19.3099
$ws.Range("E35").Value = 16.602
$ws.Range("B37").Value = 8.575700000000007
$ws.Range("B38").Value = 4.719899999999999
$ws.Range("E38").Value = 16.523
$ws.Range("E43").Value = 17.09310000000001
$ws.Range("E44").Value = 16.60919999999998
$ws.Range("E47").Value = 16.4407
$ws.Range("B51").Value = 5.712299999999998
$ws.Range("E51").Value = 17.28630000000001
$ws.Range("B52").Value = 5.3033
$ws.Range("B55").Value = 5.514599999999996
$ws.Range("E57").Value = 16.57270000000001
$ws.Range("E63").Value = 18.40970000000002
$ws.Range("A69").Value = -21.61
$ws.Range("B69").Value = 5.448199999999998
$ws.Range("B70").Value = 7.389700000000005
$ws.Range("E70").Value = 17.65770000000002
$ws.Range("A76").Value = -19.65399999999999
$ws.Range("A78").Value = -19.98859999999998
$ws.Range("B81").Value = 5.507100000000004
$ws.Range("A82").Value = -21.8783
$ws.Range("A83").Value = -21.62479999999999
$ws.Range("B83").Value = 5.556100000000003
$ws.Range("E88").Value = 16.53929999999999
$ws.Range("A93").Value = -20.73579999999998
$ws.Range("E99").Value = 16.50530000000001
$ws.Range("B102").Value = 8.380300000000009
